$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("manipulate the elections. For example, ", $true, $false, $false, $false, $false, $true, 1, $false, "manipulate the elections. For example, in Mindanao *hi Mikha could you put it here* ", 2)

Write-Output "Found: $found"
